# Code Merge Changes - 9/25/2017
#
# Updates the "TestResultExcelFilePath" value (column H, row 2) on every
# payroll-run / average-weekly-earnings worksheet so it points at the new
# shared-drive location for the statutory-scenarios automation results
# workbook, replacing the old path.

$wb = $excel.ActiveWorkbook

# Old value being replaced (for reference):
# F:\\Automation NI Reports\\HMRCTestData\Statutory scenarios for 201718\\Automation Test Result for Statutory Scenarios201718.xlsx
$newPath = "F:\\Automation_TestResults\\Payroll_Tax_StatutoryScenarios\\Automation Test Result for Statutory Scenarios201718.xlsx"

$sheetNames = @(
    "ProcessPayrolFor11WeeklyShPP",
    "ProcessPayrolFor12WeeklyShPP",
    "ProcessPayrolFor13WeeklyShPP",
    "ProcessPayrolFor14WeeklyShPP",
    "ProcessPayrolFor15WeeklyShPP",
    "ProcessPayrolFor16WeeklyShPP",
    "ProcessPayrolFor17WeeklyShPP",
    "ProcessPayrolFor18WeeklyShPP",
    "ProcessPayrolFor19WeeklyShPP",
    "AverageWeeklyEarningsTestReport",
    "ProcessPayrolFor46WeeklyShPP",
    "ProcessPayrolFor47WeeklyShPP",
    "ProcessPayrolFor48WeeklyShPP"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H2").Value = $newPath
}
